$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FFIC")
$ws.Columns("D:D").Insert()

# Row 7 (date) style copy from E7
$ws.Range("D7").NumberFormat = $ws.Range("E7").NumberFormat
$ws.Range("D7").Font.Name = $ws.Range("E7").Font.Name
$ws.Range("D7").Font.Size = $ws.Range("E7").Font.Size
$ws.Range("D7").Font.Bold = $ws.Range("E7").Font.Bold

# Row 8 (data) style copy from E8
$ws.Range("D8").NumberFormat = $ws.Range("E8").NumberFormat
$ws.Range("D8").Font.Name = $ws.Range("E8").Font.Name
$ws.Range("D8").Font.Size = $ws.Range("E8").Font.Size
$ws.Range("D8").Font.Bold = $ws.Range("E8").Font.Bold
$ws.Range("D8").HorizontalAlignment = $ws.Range("E8").HorizontalAlignment

Write-Host "done"
